$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01134666666666667
$ws.Range("H2").Value = 0.03404
$ws.Range("I2").Value = 0.001209510404472147
$ws.Range("J2").Value = 0.001209510404472147
$ws.Range("M2").Value = 0.1809866666666667
$ws.Range("N2").Value = 0.54296
$ws.Range("O2").Value = 0.03987407676082905
$ws.Range("P2").Value = 0.03987407676082905
$ws.Range("Q2").Value = 0.002053595377777777
$ws.Range("R2").Value = 0.0184823584
$ws.Range("S2").Value = 0.00004822811071094378
$ws.Range("T2").Value = 0.00004822811071094378

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01134666666666667
$ws.Range("H3").Value = 0.03404
$ws.Range("I3").Value = 0.001209510404472147
$ws.Range("J3").Value = 0.001209510404472147
$ws.Range("O3").Value = 0.1057193993302571
$ws.Range("P3").Value = 0.1057193993302571
$ws.Range("Q3").Value = 0.005444762297777778
$ws.Range("R3").Value = 0.04900286068
$ws.Range("S3").Value = 0.0001278687134444917
$ws.Range("T3").Value = 0.0001278687134444917

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01134666666666667
$ws.Range("H4").Value = 0.03404
$ws.Range("I4").Value = 0.001209510404472147
$ws.Range("J4").Value = 0.001209510404472147
$ws.Range("M4").Value = 3.878113333333333
$ws.Range("N4").Value = 11.63434
$ws.Range("O4").Value = 0.8544065239089139
$ws.Range("P4").Value = 0.8544065239089139
$ws.Range("Q4").Value = 0.04400365928888889
$ws.Range("R4").Value = 0.3960329336
$ws.Range("S4").Value = 0.001033413580316711
$ws.Range("T4").Value = 0.001033413580316711

$ws.Range("I5").Value = 0.8865539289740954
$ws.Range("J5").Value = 0.8865539289740952
$ws.Range("M5").Value = 0.1809866666666667
$ws.Range("N5").Value = 0.54296
$ws.Range("O5").Value = 0.03987407676082905
$ws.Range("P5").Value = 0.03987407676082905
$ws.Range("Q5").Value = 1.505256212728889
$ws.Range("R5").Value = 13.54730591456
$ws.Range("S5").Value = 0.03535051941652767
$ws.Range("T5").Value = 0.03535051941652766

$ws.Range("I6").Value = 0.8865539289740954
$ws.Range("J6").Value = 0.8865539289740952
$ws.Range("O6").Value = 0.1057193993302571
$ws.Range("P6").Value = 0.1057193993302571
$ws.Range("S6").Value = 0.09372594884502078
$ws.Range("T6").Value = 0.09372594884502077

$ws.Range("I7").Value = 0.8865539289740954
$ws.Range("J7").Value = 0.8865539289740952
$ws.Range("M7").Value = 3.878113333333333
$ws.Range("N7").Value = 11.63434
$ws.Range("O7").Value = 0.8544065239089139
$ws.Range("P7").Value = 0.8544065239089139
$ws.Range("Q7").Value = 32.25405658980445
$ws.Range("R7").Value = 290.28650930824
$ws.Range("S7").Value = 0.757477460712547
$ws.Range("T7").Value = 0.7574774607125468

$ws.Range("G8").Value = 1.052914333333334
$ws.Range("H8").Value = 3.158743
$ws.Range("I8").Value = 0.1122365606214325
$ws.Range("J8").Value = 0.1122365606214325
$ws.Range("M8").Value = 0.1809866666666667
$ws.Range("N8").Value = 0.54296
$ws.Range("O8").Value = 0.03987407676082905
$ws.Range("P8").Value = 0.03987407676082905
$ws.Range("Q8").Value = 0.1905634554755556
$ws.Range("R8").Value = 1.71507109928
$ws.Range("S8").Value = 0.004475329233590444
$ws.Range("T8").Value = 0.004475329233590444

$ws.Range("G9").Value = 1.052914333333334
$ws.Range("H9").Value = 3.158743
$ws.Range("I9").Value = 0.1122365606214325
$ws.Range("J9").Value = 0.1122365606214325
$ws.Range("O9").Value = 0.1057193993302571
$ws.Range("P9").Value = 0.1057193993302571
$ws.Range("Q9").Value = 0.5052469093645556
$ws.Range("R9").Value = 4.547222184281001
$ws.Range("S9").Value = 0.01186558177179183
$ws.Range("T9").Value = 0.01186558177179183

$ws.Range("G10").Value = 1.052914333333334
$ws.Range("H10").Value = 3.158743
$ws.Range("I10").Value = 0.1122365606214325
$ws.Range("J10").Value = 0.1122365606214325
$ws.Range("M10").Value = 3.878113333333333
$ws.Range("N10").Value = 11.63434
$ws.Range("O10").Value = 0.8544065239089139
$ws.Range("P10").Value = 0.8544065239089139
$ws.Range("Q10").Value = 4.083321114957778
$ws.Range("R10").Value = 36.74989003462
$ws.Range("S10").Value = 0.09589564961605025
$ws.Range("T10").Value = 0.09589564961605025

